$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.167.69'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.862.07'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3102'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.44%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07638'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.69'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08355'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.859.18'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.187'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -2.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7075'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.95%  '
$ws.Range("E15").Value = '  -0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.164.70'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.99%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.915'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.44%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '242.96'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007801'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.114.51'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -0.94%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9993'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.867'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1584'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.93%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.31'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.329'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.401'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.265'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.42%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05141'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7960'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +9.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.913'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.44%  '
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.681'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01844'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.693'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.166.03'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -5.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.209'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +0.16%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8899'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.86'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9993'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.07'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.007.67'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5203'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.775'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.337'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000120'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4272'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.31%  '
